# Generate Report for Handoff
# Updates the "Status" cells from "In Translation" to "Ready for handoff"
# and refreshes the related "Latest ... Datetime" timestamps, widening the
# relevant status columns on each sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn (E2) and de-de (F2) status columns ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-17 16:57:33"

$wsOverview.Columns.Item(5).ColumnWidth = 16.38265482584637
$wsOverview.Columns.Item(6).ColumnWidth = 16.38265482584637

# --- zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2) ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-17 16:57:27"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.38265482584637

# --- de-de sheet: Status (C2) and Latest Handoff Datetime (H2) ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-17 16:57:33"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.38265482584637
